$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells F1:K1
$headers = @{
    "F1" = "C4.5 acc"
    "G1" = "credal-C4.5 acc"
    "H1" = "SPN acc"
    "I1" = "CSPN low"
    "J1" = "CSPN high"
    "K1" = "CSPN robust"
}
# Copy the header style (bold, centered, bordered) from A1 into F1:K1
$ws.Range("A1").Copy()
$ws.Range("F1:K1").PasteSpecial(-4122)

foreach ($addr in $headers.Keys) {
    $ws.Range($addr).Value = $headers[$addr]
}

# Updated values for B2:E6
$data = @{
    "B2" = 50.66666666666666
    "C2" = 32.56410256410257
    "D2" = 82.10256410256412
    "E2" = 71.02455826605423

    "B3" = 49.94871794871795
    "C3" = 30.46153846153846
    "D3" = 82.2051282051282
    "E3" = 70.21985197549131

    "B4" = 48.05128205128206
    "C4" = 30.1025641025641
    "D4" = 81.12820512820512
    "E4" = 67.65002635743639

    "B5" = 51.28205128205128
    "C5" = 31.58974358974359
    "D5" = 80.87179487179488
    "E5" = 69.44227506775074

    "B6" = 50.15384615384615
    "C6" = 30.56410256410257
    "D6" = 81.7948717948718
    "E6" = 68.06291273582693

    "F2" = 61.8974358974359
    "G2" = 60.82051282051282
    "H2" = 57.12820512820512
    "I2" = 57.12820512820512
    "J2" = 57.12820512820512
    "K2" = 57.12820512820512

    "F3" = 59.28205128205128
    "G3" = 59.17948717948718
    "H3" = 55.69230769230769
    "I3" = 56.92307692307692
    "J3" = 57.53846153846154
    "K3" = 57.42452275920018

    "F4" = 53.17948717948718
    "G4" = 57.38461538461538
    "H4" = 58.87179487179488
    "I4" = 58.56410256410256
    "J4" = 59.33333333333334
    "K4" = 59.33182677537516

    "F5" = 46
    "G5" = 52.66666666666666
    "H5" = 54.92307692307692
    "I5" = 56.1025641025641
    "J5" = 57.28205128205128
    "K5" = 57.07297492680557

    "F6" = 43.69230769230769
    "G6" = 53.33333333333334
    "H6" = 54.82051282051282
    "I6" = 55.69230769230769
    "J6" = 58.2051282051282
    "K6" = 57.551457788067
}
foreach ($addr in $data.Keys) {
    $ws.Range($addr).Value = $data[$addr]
}
